# Auto-generated script to apply 2025-08-23 YTD crime data updates
# across Citywide Totals, By Neighborhood, and individual neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4345
$ws.Range("L3").Value = 4597
$ws.Range("L4").Value = 1138
$ws.Range("L5").Value = 264
$ws.Range("L6").Value = 3970
$ws.Range("L7").Value = 14314

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 168

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 277
$ws.Range("L3").Value = 317
$ws.Range("L7").Value = 952

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 226
$ws.Range("L6").Value = 203
$ws.Range("L7").Value = 665

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L4").Value = 11
$ws.Range("L7").Value = 191

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 155
$ws.Range("L3").Value = 183
$ws.Range("L6").Value = 148
$ws.Range("L7").Value = 532

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L6").Value = 78
$ws.Range("L7").Value = 272

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 64
$ws.Range("L6").Value = 50
$ws.Range("L7").Value = 241

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 68

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L4").Value = 56
$ws.Range("L8").Value = 952
$ws.Range("L11").Value = 230
$ws.Range("L20").Value = 363
$ws.Range("L23").Value = 149
$ws.Range("L29").Value = 802
$ws.Range("L30").Value = 68
$ws.Range("L31").Value = 140
$ws.Range("L33").Value = 665
$ws.Range("L36").Value = 180
$ws.Range("L37").Value = 532
$ws.Range("L44").Value = 102
$ws.Range("L45").Value = 24
$ws.Range("L47").Value = 104
$ws.Range("L52").Value = 290
$ws.Range("L53").Value = 168
$ws.Range("L54").Value = 297
$ws.Range("L60").Value = 89
$ws.Range("L62").Value = 9
$ws.Range("L63").Value = 44
$ws.Range("L65").Value = 272
$ws.Range("L73").Value = 115
$ws.Range("L76").Value = 218
$ws.Range("L78").Value = 185
$ws.Range("L80").Value = 46
$ws.Range("L85").Value = 742
$ws.Range("L91").Value = 199
$ws.Range("L93").Value = 76
$ws.Range("L94").Value = 179
$ws.Range("L95").Value = 191
$ws.Range("L98").Value = 80
$ws.Range("L99").Value = 241
$ws.Range("L101").Value = 14314

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 54
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 141
$ws.Range("L4").Value = 35

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L6").Value = 145
$ws.Range("L7").Value = 297

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 301
$ws.Range("L4").Value = 39
$ws.Range("L6").Value = 208
$ws.Range("L7").Value = 802

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 39
$ws.Range("L4").Value = 28
$ws.Range("L7").Value = 218

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 148
$ws.Range("L4").Value = 33

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 185

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 149

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 71
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 199

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 110
$ws.Range("L3").Value = 116
$ws.Range("L6").Value = 99
$ws.Range("L7").Value = 363

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L6").Value = 48
$ws.Range("L7").Value = 180

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 76

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L6").Value = 73
$ws.Range("L7").Value = 179

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L4").Value = 9
$ws.Range("L6").Value = 24

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 80

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 87
$ws.Range("L3").Value = 70
$ws.Range("L6").Value = 55
$ws.Range("L7").Value = 230

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 115

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L2").Value = 29
$ws.Range("L3").Value = 31
$ws.Range("L7").Value = 89

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 225
$ws.Range("L7").Value = 742

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L2").Value = 8
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 98
$ws.Range("L7").Value = 290

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L5").Value = 1
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item("Museum Campus")
$ws.Range("L3").Value = 2
$ws.Range("L6").Value = 2
$ws.Range("L7").Value = 9
